$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "58.205.56"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +3.17%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.342.50"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +0.67%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "543.92"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +6.33%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "135.17"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +2.20%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +0.84%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "2.362.60"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +1.45%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "5.43"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +3.82%  "
$ws.Range("E12").Value = "  +1.07%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.354"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +5.30%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "2.772.95"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +1.14%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "23.61"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.33%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "58.123.48"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +3.07%  "
$ws.Range("E17").Value = "  +1.55%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "2.400.98"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +2.06%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "338.89"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +4.95%  "
$ws.Range("E20").Value = "  +1.74%  "
$ws.Range("E21").Value = "  +2.00%  "
$ws.Range("E22").Value = "  +4.14%  "
$ws.Range("E23").Value = "  -0.31%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "62.15"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +1.50%  "
$ws.Range("E25").Value = "  +3.78%  "
$ws.Range("E26").Value = "  -0.29%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.15%  "
$ws.Range("E28").Value = "  +8.93%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.76"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +6.11%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "171.41"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +2.56%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.0₃0735"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +2.73%  "
$ws.Range("E32").Value = "  +1.77%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "18.55"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +1.67%  "
$ws.Range("E34").Value = "  +16.34%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.07%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("E37").Value = "  +5.82%  "
$ws.Range("E38").Value = "  +0.12%  "
$ws.Range("E39").Value = "  +4.43%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "39.43"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +2.59%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "149.55"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.38%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.378"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +1.61%  "
$ws.Range("E43").Value = "  +2.70%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "283.33"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +1.66%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "19.40"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +7.79%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.0931"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.79%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.0505"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +2.13%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.560"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +1.27%  "
$ws.Range("E49").Value = "  +2.51%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "17.59"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +4.01%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.383"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.70%  "
